# Add a new "Code Article" (default_code) column as the new column A.
# Everything that was in column A..F shifts right to B..G, e.g. the
# "Entrepôt YourCompany" title cell moves from B1 to C1, and the header
# row / data row gain a new first column holding the product's internal
# reference code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at the left (pushes existing columns A:.. to B:..)
$ws.Columns("A:A").Insert()

# New header cell for the inserted column
$ws.Range("A2").Value = "Code Article"

# New data cell holding the product's default_code, matching the existing
# "Flipover" row (now shifted to column B on row 3)
$ws.Range("A3").Value = "FURN_9001"
